# The presentation currently uses the "Integral" design theme (ppt/theme/theme1.xml,
# bound to the one and only Slide Master). The commit swaps the presentation over to
# the default Office color palette - i.e. re-colors the theme that backs the Slide
# Master from "Integral" colors to the stock "Office Theme" colors.
#
# PowerPoint's Design/Variants gallery does this by rewriting the 12 theme colors
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) of the active theme's color scheme -
# exactly what ThemeColorScheme exposes over COM. We set each slot to the
# corresponding "Office Theme" RGB value.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function BgrInt([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: stock "Office Theme" colors (in clrScheme order).
$officeThemeColors = @(
    @{ Name = "dk1";      R = 0x00; G = 0x00; B = 0x00 },  # 000000
    @{ Name = "lt1";      R = 0xFF; G = 0xFF; B = 0xFF },  # FFFFFF
    @{ Name = "dk2";      R = 0x44; G = 0x54; B = 0x6A },  # 44546A
    @{ Name = "lt2";      R = 0xE7; G = 0xE6; B = 0xE6 },  # E7E6E6
    @{ Name = "accent1";  R = 0x5B; G = 0x9B; B = 0xD5 },  # 5B9BD5
    @{ Name = "accent2";  R = 0xED; G = 0x7D; B = 0x31 },  # ED7D31
    @{ Name = "accent3";  R = 0xA5; G = 0xA5; B = 0xA5 },  # A5A5A5
    @{ Name = "accent4";  R = 0xFF; G = 0xC0; B = 0x00 },  # FFC000
    @{ Name = "accent5";  R = 0x44; G = 0x72; B = 0xC4 },  # 4472C4
    @{ Name = "accent6";  R = 0x70; G = 0xAD; B = 0x47 },  # 70AD47
    @{ Name = "hlink";    R = 0x05; G = 0x63; B = 0xC1 },  # 0563C1
    @{ Name = "folHlink"; R = 0x95; G = 0x4F; B = 0x72 }   # 954F72
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $entry = $officeThemeColors[$i]
    $colorScheme.Item($i + 1).RGB = BgrInt $entry.R $entry.G $entry.B
}
